# Auto-generated Excel COM-interop script applying the Jenova_Profits.xlsx diff
# (currentAveragePrice / LevePrice / LeveProfit recompute across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 44648.78
$ws.Range("J28").Value = 4999.5
$ws.Range("L28").Value = 4999.5
$ws.Range("N28").Value = -5969.5

$ws.Range("H32").Value = 7649.1665
$ws.Range("J32").Value = 8179
$ws.Range("L32").Value = 8179
$ws.Range("N32").Value = -8831

$ws.Range("H40").Value = 9654.154
$ws.Range("I40").Value = 8428.857
$ws.Range("J40").Value = 11083.667
$ws.Range("K40").Value = 8428.857
$ws.Range("L40").Value = 11083.667
$ws.Range("M40").Value = -8253.857
$ws.Range("N40").Value = -11433.667

$ws.Range("H44").Value = 565555.5
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H62").Value = 6252524
$ws.Range("I62").Value = 8930484
$ws.Range("K62").Value = 8930484
$ws.Range("M62").Value = -8929860

$ws.Range("H65").Value = 6252524
$ws.Range("I65").Value = 8930484
$ws.Range("K65").Value = 44652420
$ws.Range("M65").Value = -44649300

$ws.Range("H70").Value = 201319.6
$ws.Range("J70").Value = 201319.6
$ws.Range("L70").Value = 603958.8
$ws.Range("N70").Value = -604498.8

$ws.Range("H73").Value = 201319.6
$ws.Range("J73").Value = 201319.6
$ws.Range("L73").Value = 603958.8
$ws.Range("N73").Value = -605830.8

$ws.Range("H74").Value = 7691.7856
$ws.Range("I74").Value = 4153.1816
$ws.Range("J74").Value = 20666.666
$ws.Range("K74").Value = 4153.1816
$ws.Range("L74").Value = 20666.666
$ws.Range("M74").Value = -3217.1816
$ws.Range("N74").Value = -22538.666

$ws.Range("H76").Value = 111231640
$ws.Range("I76").Value = 154253.72
$ws.Range("K76").Value = 154253.72
$ws.Range("M76").Value = -153938.72

$ws.Range("H77").Value = 7691.7856
$ws.Range("I77").Value = 4153.1816
$ws.Range("J77").Value = 20666.666
$ws.Range("K77").Value = 20765.908
$ws.Range("L77").Value = 103333.33
$ws.Range("M77").Value = -16085.908
$ws.Range("N77").Value = -112693.33

$ws.Range("H79").Value = 111231640
$ws.Range("I79").Value = 154253.72
$ws.Range("K79").Value = 154253.72
$ws.Range("M79").Value = -153161.72

$ws.Range("H112").Value = 4171.4194
$ws.Range("J112").Value = 4486.2144
$ws.Range("L112").Value = 13458.6432
$ws.Range("N112").Value = -15674.6432

$ws.Range("H132").Value = 2808.9092
$ws.Range("I132").Value = 2159.2974
$ws.Range("K132").Value = 6477.8922
$ws.Range("M132").Value = -3947.8922

$ws.Range("H137").Value = 6009.2
$ws.Range("I137").Value = 6248
$ws.Range("J137").Value = 5850
$ws.Range("K137").Value = 18744
$ws.Range("L137").Value = 17550
$ws.Range("M137").Value = -16194
$ws.Range("N137").Value = -22650

$ws.Range("H138").Value = 5628.6084
$ws.Range("I138").Value = 2646.2222
$ws.Range("J138").Value = 6778.957
$ws.Range("K138").Value = 7938.6666
$ws.Range("L138").Value = 20336.871
$ws.Range("M138").Value = -2798.6666
$ws.Range("N138").Value = -30616.871

$ws.Range("H141").Value = 3032.5334
$ws.Range("I141").Value = 2300.9167
$ws.Range("K141").Value = 6902.750100000001
$ws.Range("M141").Value = -1722.750100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 10
$ws.Range("J4").Value = 10
$ws.Range("L4").Value = 10
$ws.Range("N4").Value = -242

$ws.Range("H23").Value = 1719715.1
$ws.Range("I23").Value = 4000002
$ws.Range("K23").Value = 4000002
$ws.Range("M23").Value = -3999743

$ws.Range("H32").Value = 4802.2793
$ws.Range("I32").Value = 3605.303
$ws.Range("K32").Value = 3605.303
$ws.Range("M32").Value = -3318.303

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H63").Value = 7710.3
$ws.Range("I63").Value = 5314.7144
$ws.Range("J63").Value = 13300
$ws.Range("K63").Value = 5314.7144
$ws.Range("L63").Value = 13300
$ws.Range("M63").Value = -4628.7144
$ws.Range("N63").Value = -14672

$ws.Range("H66").Value = 7710.3
$ws.Range("I66").Value = 5314.7144
$ws.Range("K66").Value = 26573.572
$ws.Range("L66").Value = 66500
$ws.Range("M66").Value = -23141.572
$ws.Range("N66").Value = -73364

$ws.Range("H122").Value = 2885.6667
$ws.Range("J122").Value = 3767.0435
$ws.Range("L122").Value = 11301.1305
$ws.Range("N122").Value = -16201.1305

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 54306.156
$ws.Range("I105").Value = 72566.57000000001
$ws.Range("J105").Value = 3177
$ws.Range("K105").Value = 72566.57000000001
$ws.Range("L105").Value = 3177
$ws.Range("M105").Value = -70819.57000000001
$ws.Range("N105").Value = -6671

$ws.Range("H132").Value = 49999.832
$ws.Range("J132").Value = 49999.832
$ws.Range("L132").Value = 49999.832
$ws.Range("N132").Value = -60119.832

$ws.Range("H134").Value = 25096.646
$ws.Range("I134").Value = 3108.5151
$ws.Range("K134").Value = 9325.5453
$ws.Range("M134").Value = -6790.5453

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5678.5713
$ws.Range("I62").Value = 6372
$ws.Range("K62").Value = 6372
$ws.Range("M62").Value = -5748

$ws.Range("H65").Value = 5678.5713
$ws.Range("I65").Value = 6372
$ws.Range("K65").Value = 31860
$ws.Range("M65").Value = -28740

$ws.Range("H132").Value = 3221.2285
$ws.Range("J132").Value = 5655.5835
$ws.Range("L132").Value = 16966.7505
$ws.Range("N132").Value = -22026.7505

$ws.Range("H134").Value = 221343.2
$ws.Range("I134").Value = 2800.16
$ws.Range("K134").Value = 8400.48
$ws.Range("M134").Value = -5865.48

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1263
$ws.Range("I36").Value = 828
$ws.Range("J36").Value = 3003
$ws.Range("K36").Value = 2484
$ws.Range("L36").Value = 9009
$ws.Range("M36").Value = -2315
$ws.Range("N36").Value = -9347

$ws.Range("H92").Value = 526891.6
$ws.Range("I92").Value = 1667081.9
$ws.Range("J92").Value = 650
$ws.Range("K92").Value = 5001245.699999999
$ws.Range("L92").Value = 1950
$ws.Range("M92").Value = -4999997.699999999
$ws.Range("N92").Value = -4446

$ws.Range("H97").Value = 761.375
$ws.Range("I97").Value = 1208
$ws.Range("J97").Value = 493.4
$ws.Range("K97").Value = 3624
$ws.Range("L97").Value = 1480.2
$ws.Range("M97").Value = -3128
$ws.Range("N97").Value = -2472.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 81.42856999999999
$ws.Range("I2").Value = 84.90000000000001
$ws.Range("K2").Value = 84.90000000000001
$ws.Range("M2").Value = 28.09999999999999

$ws.Range("H7").Value = 6143213.5
$ws.Range("J7").Value = 2181817.2
$ws.Range("L7").Value = 2181817.2
$ws.Range("N7").Value = -2182041.2

$ws.Range("H8").Value = 6143213.5
$ws.Range("J8").Value = 2181817.2
$ws.Range("L8").Value = 2181817.2
$ws.Range("N8").Value = -2182095.2

$ws.Range("H80").Value = 1056378
$ws.Range("I80").Value = 628960.7
$ws.Range("K80").Value = 628960.7
$ws.Range("M80").Value = -627962.7

$ws.Range("H83").Value = 1056378
$ws.Range("I83").Value = 628960.7
$ws.Range("K83").Value = 3144803.5
$ws.Range("M83").Value = -3139811.5

$ws.Range("H132").Value = 246922.28
$ws.Range("I132").Value = 252215.08
$ws.Range("J132").Value = 204580
$ws.Range("K132").Value = 756645.24
$ws.Range("L132").Value = 613740
$ws.Range("M132").Value = -754115.24
$ws.Range("N132").Value = -618800

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 94777.37
$ws.Range("I68").Value = 3475.125
$ws.Range("K68").Value = 3475.125
$ws.Range("M68").Value = -2726.125

$ws.Range("H71").Value = 94777.37
$ws.Range("I71").Value = 3475.125
$ws.Range("K71").Value = 17375.625
$ws.Range("M71").Value = -13631.625

$ws.Range("H132").Value = 6125.25
$ws.Range("I132").Value = 5375
$ws.Range("K132").Value = 16125
$ws.Range("M132").Value = -13595

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 32121.092
$ws.Range("I2").Value = 32121.092
$ws.Range("K2").Value = 32121.092
$ws.Range("M2").Value = -32009.092

$ws.Range("H113").Value = 2159.389
$ws.Range("I113").Value = 1615.0834
$ws.Range("J113").Value = 3248
$ws.Range("K113").Value = 4845.2502
$ws.Range("L113").Value = 9744
$ws.Range("M113").Value = -2675.2502
$ws.Range("N113").Value = -14084

$ws.Range("H136").Value = 95165.03999999999
$ws.Range("I136").Value = 22184.7
$ws.Range("K136").Value = 66554.10000000001
$ws.Range("M136").Value = -64004.10000000001
